$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1219.956671505592
$ws.Range("D2").Value = 4652.939829936181
